# Added stats reporting to excel_ui.
# Update sample G column stats and switch the active sheet/selection
# from "Beads" to "Samples" (on cell K11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Samples")

# Updated stats values in column G
$ws.Range("G8").Value = 0.2
$ws.Range("G9").Value = 0.2
$ws.Range("G10").Value = 0.25

# Make "Samples" the active/selected sheet and select cell K11
$ws.Activate()
$ws.Range("K11").Select()
